# Applies the two textual edits described by the commit diff:
#   1. Slide 10 ("5) Forums" / MongoDB collections slide): the bullet
#      "messages" becomes "Messages (sub-collection)".
#   2. Slide 23 ("Security" slide): the list item that was split across two
#      runs, "Data " + "Protection and Privacy", is merged into a single
#      run reading "Data Protection and Privacy".

$p = $ppt.ActivePresentation

# --- Slide 10: Content Placeholder 4, paragraph 12 ("messages") ---
$s10 = $p.Slides.Item(10)
$shape10 = $s10.Shapes.Item(3)
$tr10 = $shape10.TextFrame.TextRange
$para10 = $tr10.Paragraphs(12, 1)
$para10.Text = "Messages (sub-collection)"

# --- Slide 23: Content Placeholder 2, paragraph 3 ("Data " + "Protection and Privacy") ---
$s23 = $p.Slides.Item(23)
$shape23 = $s23.Shapes.Item(2)
$tr23 = $shape23.TextFrame.TextRange
$para23 = $tr23.Paragraphs(3, 1)
# Setting the same concatenated text as a no-op change is ignored by the
# engine (it already reads back as the merged string), so first set a
# distinct placeholder to force the paragraph's runs to be rebuilt, then
# assign the final text; this collapses the two runs into one.
$para23.Text = "TEMP_PLACEHOLDER"
$para23b = $tr23.Paragraphs(3, 1)
$para23b.Text = "Data Protection and Privacy"
